# Auto-generated edit script: refresh market-price derived columns (H-N)
# across the Behemoth_Profits workbook sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ===== Sheet ALC =====
$ws = $wb.Worksheets.Item("ALC")
# row 17
$ws.Range("H17").Value = 568.4545000000001
$ws.Range("J17").Value = 568.4545000000001
$ws.Range("L17").Value = 1705.3635
$ws.Range("N17").Value = -2041.3635
# row 28
$ws.Range("H28").Value = 962.5
$ws.Range("I28").Value = 642
$ws.Range("J28").Value = 1763.75
$ws.Range("K28").Value = 642
$ws.Range("L28").Value = 1763.75
$ws.Range("M28").Value = -157
$ws.Range("N28").Value = -2733.75
# row 33
$ws.Range("H33").Value = 2139.8108
$ws.Range("I33").Value = 2134.0908
$ws.Range("J33").Value = 2187
$ws.Range("K33").Value = 2134.0908
$ws.Range("L33").Value = 2187
$ws.Range("M33").Value = -1905.0908
$ws.Range("N33").Value = -2645
# row 55
$ws.Range("H55").Value = 1298.6471
$ws.Range("I55").Value = 132.33333
$ws.Range("J55").Value = 4097.8
$ws.Range("K55").Value = 132.33333
$ws.Range("L55").Value = 4097.8
$ws.Range("M55").Value = 81.66667000000001
$ws.Range("N55").Value = -4525.8
# row 64
$ws.Range("H64").Value = 5022.1665
# row 67
$ws.Range("H67").Value = 5022.1665
# row 126
$ws.Range("H126").Value = 123000
$ws.Range("J126").Value = 123000
$ws.Range("L126").Value = 123000
$ws.Range("N126").Value = -132880
# row 130
$ws.Range("H130").Value = 71999.664
$ws.Range("J130").Value = 71999.664
$ws.Range("L130").Value = 71999.664
$ws.Range("N130").Value = -82039.664
# row 132
$ws.Range("H132").Value = 2376.8667
$ws.Range("I132").Value = 1885.8462
$ws.Range("J132").Value = 5568.5
$ws.Range("K132").Value = 5657.5386
$ws.Range("L132").Value = 16705.5
$ws.Range("M132").Value = -3127.5386
$ws.Range("N132").Value = -21765.5
# row 138
$ws.Range("H138").Value = 1820.6666
$ws.Range("I138").Value = 1032.5238
$ws.Range("K138").Value = 3097.5714
$ws.Range("M138").Value = 2042.4286

# ===== Sheet ARM =====
$ws = $wb.Worksheets.Item("ARM")
# row 17
$ws.Range("H17").Value = 10000
$ws.Range("I17").Value = 10000
$ws.Range("K17").Value = 10000
$ws.Range("M17").Value = -9827
# row 22
$ws.Range("H22").Value = 1016
$ws.Range("I22").Value = 1016
$ws.Range("K22").Value = 1016
$ws.Range("M22").Value = -717
# row 125
$ws.Range("H125").Value = 70715
$ws.Range("J125").Value = 70715
$ws.Range("L125").Value = 70715
$ws.Range("N125").Value = -80555
# row 139
$ws.Range("H139").Value = 59000
$ws.Range("I139").Value = 59000
$ws.Range("K139").Value = 59000
$ws.Range("M139").Value = -53860

# ===== Sheet CRP =====
$ws = $wb.Worksheets.Item("CRP")
# row 7
$ws.Range("H7").Value = 2924.6667
$ws.Range("I7").Value = 180.22223
$ws.Range("J7").Value = 11158
$ws.Range("K7").Value = 180.22223
$ws.Range("L7").Value = 11158
$ws.Range("M7").Value = -67.22223
$ws.Range("N7").Value = -11384
# row 112
$ws.Range("H112").Value = 64308.6
$ws.Range("J112").Value = 64308.6
$ws.Range("L112").Value = 64308.6
$ws.Range("N112").Value = -67262.60000000001
# row 134
$ws.Range("H134").Value = 504771.34
$ws.Range("I134").Value = 716037
$ws.Range("K134").Value = 2148111
$ws.Range("M134").Value = -2145576

# ===== Sheet CUL =====
$ws = $wb.Worksheets.Item("CUL")
# row 11
$ws.Range("H11").Value = 2341.818
$ws.Range("I11").Value = 2396.1904
$ws.Range("K11").Value = 7188.5712
$ws.Range("M11").Value = -7048.5712
# row 37
$ws.Range("H37").Value = 84991
$ws.Range("J37").Value = 84991
$ws.Range("L37").Value = 254973
$ws.Range("N37").Value = -255197
# row 41
$ws.Range("H41").Value = 1083.1666
$ws.Range("I41").Value = 966.6667
$ws.Range("J41").Value = 1199.6666
$ws.Range("K41").Value = 2900.0001
$ws.Range("L41").Value = 3598.9998
$ws.Range("M41").Value = -2562.0001
$ws.Range("N41").Value = -4274.9998
# row 52
$ws.Range("H52").Value = 8516.25
$ws.Range("J52").Value = 8516.25
$ws.Range("L52").Value = 25548.75
$ws.Range("N52").Value = -26080.75
# row 93
$ws.Range("H93").Value = 5141.4443
$ws.Range("I93").Value = 2000
$ws.Range("J93").Value = 5534.125
$ws.Range("K93").Value = 6000
$ws.Range("L93").Value = 16602.375
$ws.Range("M93").Value = -4128
$ws.Range("N93").Value = -20346.375
# row 114
$ws.Range("H114").Value = 1211.6666
$ws.Range("I114").Value = 867.75
$ws.Range("K114").Value = 2603.25
$ws.Range("M114").Value = 650.75
# row 115
$ws.Range("H115").Value = 49354
$ws.Range("I115").Value = 8000
$ws.Range("J115").Value = 70031
$ws.Range("K115").Value = 24000
$ws.Range("L115").Value = 210093
$ws.Range("M115").Value = -22825
$ws.Range("N115").Value = -212443
# row 139
$ws.Range("H139").Value = 2600
$ws.Range("I139").Value = 3500
$ws.Range("K139").Value = 10500
$ws.Range("M139").Value = -5360
# row 140
$ws.Range("H140").Value = 275254.47
$ws.Range("I140").Value = 302364.9
$ws.Range("K140").Value = 907094.7000000001
$ws.Range("M140").Value = -901914.7000000001

# ===== Sheet GSM =====
$ws = $wb.Worksheets.Item("GSM")
# row 93
$ws.Range("H93").Value = 60000
$ws.Range("J93").Value = 60000
$ws.Range("L93").Value = 60000
$ws.Range("N93").Value = -63744
# row 102
$ws.Range("H102").Value = 4187.4595
$ws.Range("I102").Value = 3038.7036
$ws.Range("J102").Value = 7289.1
$ws.Range("K102").Value = 3038.7036
$ws.Range("L102").Value = 7289.1
$ws.Range("M102").Value = -1416.7036
$ws.Range("N102").Value = -10533.1

# ===== Sheet LTW =====
$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 2499.5557
$ws.Range("I22").Value = 2565.2104
$ws.Range("J22").Value = 2343.625
$ws.Range("K22").Value = 2565.2104
$ws.Range("L22").Value = 2343.625
$ws.Range("M22").Value = -2270.2104
$ws.Range("N22").Value = -2933.625
# row 27
$ws.Range("H27").Value = 2499.5557
$ws.Range("I27").Value = 2565.2104
$ws.Range("J27").Value = 2343.625
$ws.Range("K27").Value = 2565.2104
$ws.Range("L27").Value = 2343.625
$ws.Range("M27").Value = -2458.2104
$ws.Range("N27").Value = -2557.625
# row 55
$ws.Range("H55").Value = 50000228
$ws.Range("I55").Value = 66666896
$ws.Range("K55").Value = 66666896
$ws.Range("M55").Value = -66666723
# row 61
$ws.Range("H61").Value = 959.8
$ws.Range("I61").Value = 959.8
$ws.Range("K61").Value = 959.8
$ws.Range("M61").Value = -757.8
# row 102
$ws.Range("H102").Value = 79999.5
$ws.Range("I102").Value = 79999
$ws.Range("J102").Value = 80000
$ws.Range("K102").Value = 79999
$ws.Range("L102").Value = 80000
$ws.Range("M102").Value = -76754
$ws.Range("N102").Value = -86490
# row 113
$ws.Range("H113").Value = 959.8
$ws.Range("I113").Value = 959.8
$ws.Range("K113").Value = 959.8
$ws.Range("M113").Value = 1210.2

# ===== Sheet WVR =====
$ws = $wb.Worksheets.Item("WVR")
# row 43
$ws.Range("H43").Value = 80000
$ws.Range("I43").Value = 80000
$ws.Range("K43").Value = 80000
$ws.Range("M43").Value = -79851
# row 96
$ws.Range("H96").Value = 3999.25
$ws.Range("I96").Value = 1999
$ws.Range("K96").Value = 1999
$ws.Range("M96").Value = -626
# row 107
$ws.Range("H107").Value = 27778958
$ws.Range("J107").Value = 1397.8
$ws.Range("L107").Value = 4193.4
$ws.Range("N107").Value = -8033.4
# row 132
$ws.Range("H132").Value = 9857.166999999999
$ws.Range("I132").Value = 1661.7368
$ws.Range("K132").Value = 4985.2104
$ws.Range("M132").Value = -2455.2104
